# Rectification du mouvement des pièces
# Fill in the previously-empty rows 71-74 of the "Journal de travail"
# (Tableau1) on sheet "Feuil1": a missing end-time on row 71, and three
# brand-new entries on rows 72-74.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Rows 72-74 are brand new table rows: copy the formatting from row 71
# first (date / time / text styles) so the new cells pick up the same
# number formats as the rest of the table instead of plain defaults.
$ws.Range("B71:G71").Copy()
$ws.Range("B72:G72").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B73:G73").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B74:G74").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 71 : the Implémentation entry was missing its "Fin" (end time) - 14:00
$ws.Range("D71").Value = 0.58333333333333337

# Row 72 : new entry - Implémentation / 14:00 -> 14:30
$ws.Range("B72").Value = 45433
$ws.Range("C72").Value = 0.58333333333333337
$ws.Range("D72").Value = 0.60416666666666663
$ws.Range("F72").Value = "Implémentation"
$ws.Range("G72").Value = "Implémentation de la prévision du déplacement des pièces"

# Row 73 : new entry - Documentation / 14:30 -> 15:05
$ws.Range("B73").Value = 45433
$ws.Range("C73").Value = 0.60416666666666663
$ws.Range("D73").Value = 0.62847222222222221
$ws.Range("F73").Value = "Documentation"
$ws.Range("G73").Value = "Rédaction de la génération des mouvements des pièces"

# Row 74 : new entry - Documentation / 15:20 -> 15:55
$ws.Range("B74").Value = 45433
$ws.Range("C74").Value = 0.63888888888888895
$ws.Range("D74").Value = 0.66319444444444442
$ws.Range("F74").Value = "Documentation"
$ws.Range("G74").Value = "Rédaction de la génération des mouvements des pièces"

# Move the active selection to match the author's final cursor position
$ws.Range("H70").Select()
